$d = $word.ActiveDocument

$replacements = @(
    @("314÷4=", "652÷7="),
    @("370÷8=", "949÷4="),
    @("146÷2=", "472÷9="),
    @("132÷6=", "695÷7="),
    @("699÷8=", "999÷5="),
    @("775÷7=", "514÷2="),
    @("711÷4=", "630÷5="),
    @("979÷8=", "481÷7="),
    @("540÷6=", "464÷3="),
    @("718÷3=", "117÷7="),
    @("501÷3=", "769÷5="),
    @("881÷9=", "711÷2="),
    @("361÷6=", "591÷8="),
    @("433÷8=", "895÷5="),
    @("568÷8=", "884÷9="),
    @("443÷9=", "838÷5="),
    @("166÷5=", "940÷2="),
    @("559÷7=", "900÷7="),
    @("452÷4=", "696÷7="),
    @("755÷5=", "971÷8="),
    @("354÷7=", "284÷7="),
    @("856÷2=", "965÷3="),
    @("580÷8=", "679÷5="),
    @("505÷2=", "900÷2="),
    @("611÷7=", "901÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
